$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The demo invoices' relative day-offsets (date_invoice / date columns) were
# recomputed against a different "today" reference, shrinking every
# "-62" placeholder (and the few rows that already diverged from it) down
# to the new, smaller offsets.

# Rows 4-11 (out_invoice/out_refund demo rows) all shared the same
# "-62" placeholder in column F; it becomes "-35" everywhere it was used.
$ws.Range("F4").Value = "-35"
$ws.Range("F5").Value = "-35"
$ws.Range("F6").Value = "-35"
$ws.Range("F7").Value = "-35"
$ws.Range("F8").Value = "-35"
$ws.Range("F9").Value = "-35"
$ws.Range("F10").Value = "-35"
$ws.Range("F11").Value = "-35"

# Rows 13-19 (in_invoice/in_refund purchase demo rows): both date_invoice
# (F) and date (G) columns get new, closer-to-zero offsets.
$ws.Range("F13").Value = "-39"
$ws.Range("G13").Value = "-35"

$ws.Range("F14").Value = "-38"
$ws.Range("G14").Value = "-35"

$ws.Range("F15").Value = "-37"
$ws.Range("G15").Value = "-35"

$ws.Range("F16").Value = "-36"
$ws.Range("G16").Value = "-35"

$ws.Range("F17").Value = "-35"
$ws.Range("G17").Value = "-35"

$ws.Range("F18").Value = "-35"
$ws.Range("G18").Value = "-35"

$ws.Range("F19").Value = "-42"
$ws.Range("G19").Value = "-35"

# Move the active selection from F19:G19 to F10, matching the cursor
# position left by the author after editing.
$ws.Range("F10").Select()
